# v2p14. Compatible with MF-Swift v2212, updated hardpoints.
#
# Updates the damper hardpoint figures (sTop / sBottom rows, columns F:H)
# on both sheets, gives the sheet1 "x" hardpoint column (F) a finer
# (three-decimal) display format, and leaves the workbook focused on
# sheet1 with the edited range selected (sheet2 keeps its own last
# selection).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Sedan_HambaLG_f
$ws2 = $wb.Worksheets.Item(2)   # Sedan_HambaLG_r

# ---------------------------------------------------------------------
# Sheet1 (Sedan_HambaLG_f) - sTop (row 7) / sBottom (row 8) hardpoints
# ---------------------------------------------------------------------
$ws1.Range("F7").Value = -0.002655714285714287
$ws1.Range("G7").Value = 0.62
$ws1.Range("H7").Value = 0.65

$ws1.Range("F8").Value = 0.05516642857142858
$ws1.Range("G8").Value = 0.85
$ws1.Range("H8").Value = 0.19

# Column F (x) gets a dedicated 3-decimal number format on this sheet.
$ws1.Range("F7").NumberFormat = "0.000"
$ws1.Range("F8").NumberFormat = "0.000"
$ws1.Range("G7:H8").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# Sheet2 (Sedan_HambaLG_r) - sTop (row 7) / sBottom (row 8) hardpoints
# ---------------------------------------------------------------------
$ws2.Range("F7").Value = 0.002655714285714287
$ws2.Range("G7").Value = 0.62
$ws2.Range("H7").Value = 0.65

$ws2.Range("F8").Value = -0.05516642857142858
$ws2.Range("G8").Value = 0.85
$ws2.Range("H8").Value = 0.19

$ws2.Range("F7:H8").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# Tab colors - switch both tabs from the pale "theme 7" accent to the
# darker "theme 8" accent used in the new revision.
# ---------------------------------------------------------------------
$ws1.Tab.Color = 10515524
$ws2.Tab.Color = 10515524

# ---------------------------------------------------------------------
# Selection / active sheet - sheet1 becomes the active tab with the
# edited F7:H8 block selected; sheet2 keeps an independent selection.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("E15").Select()

$ws1.Activate()
$ws1.Range("F7:H8").Select()
